# Opioid Use Disorder Treatment Medical Billing Codes workbook:
# wrap each CPT/HCPCS billing code on Sheet2 in single quotes inside the
# per-row CONCAT helper formula (column C), which ripples into the
# aggregate bracketed list built in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Every row's helper cell in column C wraps the CPT/HCPCS code (column B)
# in single quotes before the trailing comma, e.g. 90832, -> '90832',
$lastRow = 140
for ($r = 1; $r -le $lastRow; $r++) {
    $bRef = "B" + $r
    $ws.Cells.Item($r, 3).Formula = '=_xlfn.CONCAT("''",' + $bRef + ',"''",",")'
}

# Leftover literal paste of the rebuilt aggregate list that the author left
# sitting in E4 while reviewing the new quoted output.
$ws.Range("E4").Value = "['90832','90834','G1028','G2067','G2068','G2069','G2070','G2071','G2072','G2073','G2074','G2075','G2076','G2077','G2078','G2079','G2080','G2081','G2086','G2087','G2088','G2215','G2216','H0020','H0001','H0002','H0003','H0004','H0005','H0006','H0007','H0008','H0009','H0010','H0011','H0012','H0013','H0014','H0015','H0016','H0017','H0018','H0019','H0021','H0022','H0023','H0024','H0025','H0026','H0027','H0028','H0029','H0030','H0031','H0032','H0033','H0034','H0035','H0036','H0037','H0038','H0039','H0040','H0041','H0042','H0043','H0044','H0045','H0046','H0047','H0048','H0049','H0050','H1000','H1001','H1002','H1003','H1004','H1005','H1010','H1011','H2000','H2001','H2010','H2011','H2012','H2013','H2014','H2015','H2016','H2017','H2018','H2019','H2020','H2021','H2022','H2023','H2024','H2025','H2026','H2027','H2028','H2029','H2030','H2031','H2032','H2033','H2034','H2035','H2036','H2037','J2315','T1002','T1006','T1007','T1012','T1016','T2048','90791','90837','90839','90840','90845','90846','90847','90853','90801','90802','90804','90806','90808','90810','90812','90816','90818','90821','90823','90826','90828','90857',]"

# Restore the worksheet selection the author ended up with.
$ws.Activate()
$ws.Range("G10").Select()
